$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generators")

$ws.Range("H2").Formula = "=-G2"
$ws.Range("H3:H54").Formula = "=-G3"

$ws.Range("H2:H54").Select()
